$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    ,@(163, '2026-02-14 12:36:42', '237681102046', 'FRANCOISE NKENFACK NKENGMO', 397988)
    ,@(164, '2026-02-14 14:16:37', '237682368679', 'MFS SIM AA 2', 322045)
    ,@(165, '2026-02-14 05:26:00', '237683360459', 'LUCIE MAJOLIE LELE NKANKEU', 574)
    ,@(166, '2026-02-14 13:53:18', '237652899422', 'NOUMOU epouse SAGNON MARCELINE LA NEGRESSE', 155671)
    ,@(167, '2026-02-14 15:49:46', '237670904526', 'MFS SIM PROVISOIRE 20', 0)
    ,@(168, '2026-02-14 13:22:51', '237671105116', 'MFS  AM FACE HÔPITAL GÉNÉRAL', 18)
    ,@(169, '2026-02-14 18:02:31', '237672916354', 'MAMADOU DIAN BAH LENA GLOBAL', 68564)
    ,@(170, '2026-02-14 03:29:45', '237672920086', 'NAMY NGOKO CLARISSE ROSE VERTINE KAMILAH CONNECTION', 6730)
    ,@(171, '2026-02-14 14:54:28', '237674000053', 'FRED JUNIOR ZOK EDOU', 138362)
    ,@(172, '2026-02-14 14:39:27', '237674841555', 'BEATRICE TCHAMTIEU EPSE NGAMENI', 192473)
    ,@(173, '2026-02-14 16:18:06', '237674899678', 'VIVIANE MADJUIMEKEM FOMEKONG', 82338)
    ,@(174, '2026-02-14 10:00:49', '237676439452', 'SAGNOU BRINDA JOSELINE _DIGITAL BUSINESS SARL', 229)
    ,@(175, '2026-02-14 17:08:00', '237676695935', 'SIMON PIERRE AKOA', 132533)
    ,@(176, '2026-02-14 15:40:20', '237677745809', 'FRANKLIN MUA ZUO', 228443)
    ,@(177, '2026-02-14 17:04:40', '237679127464', 'DJUFFO TSOATA MARIE NOEL KAMILAH CONNECTION GROUP', 33636)
    ,@(178, '2026-02-14 12:30:08', '237679422291', 'ETS LE CONTENT 32', 100022)
    ,@(179, '2026-02-14 12:08:43', '237651433330', 'NGUIAZONG DORIANE LAURE KAMILAH CONNECTION GROUP', 83422)
    ,@(180, '2026-02-14 15:29:04', '237654168696', 'DZEUMAZONG FLORENCE ETS MOBILE FINANCIAL SERVICES MFS', 1893)
    ,@(181, '2026-02-14 17:46:11', '237670799877', 'ETS MOBILE FINANCIAL SERVICES MFS MENANDJIO HORTENSE BIENVENUE', 642657)
    ,@(182, '2026-02-14 15:05:50', '237671351291', 'MFS LTDLA CBOX R3 MOUTHIEU JOSETTE CHANCELINE', 324495)
    ,@(183, '2026-02-14 14:28:09', '237671378136', 'KOUBINOM DIPITA SARIETTE CRISTELLE ETS MOBILE FINANCIAL SERVICES MFS', 4426)
    ,@(184, '2026-02-14 17:36:17', '237671605749', 'ETS TCHATCHOUANG PAUL  ETP LTDLA CBOX RO MEGAPTCHE VICTORINE', 436933)
    ,@(185, '2026-02-14 15:17:58', '237671615641', 'BEGO FOGUE CHRISTELLE KAMILAH CONNECTION GROUP', 21212)
    ,@(186, '2026-02-14 16:07:42', '237673739931', 'MOFFO GERMAIN SPECTRUM SPECTRUM', 45501)
)

foreach ($row in $rows) {
    $r = $row[0]
    $dateVal = $row[1]
    $numberVal = $row[2]
    $nameVal = $row[3]
    $balanceVal = $row[4]

    # Columns A-C are text (dates stored as literal strings, phone numbers as text).
    # Force text type via NumberFormat="@" then strip the format so the cell
    # ends up with no explicit style, matching a plain text write.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $dateVal
    $cellA.ClearFormats()

    $cellB = $ws.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $numberVal
    $cellB.ClearFormats()

    $cellC = $ws.Cells.Item($r, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $nameVal
    $cellC.ClearFormats()

    # Column D is numeric.
    $ws.Cells.Item($r, 4).Value = $balanceVal
}

